# Apply the "new version with timestamp" update to the daily sale shortage report.
#
# Logical change: the item that used to appear in row 11
# ("TETANUS ANTITOXIN 1500 I U 10 AMPOULES") was removed and replaced by a new
# item ("MOBITIL 15MG/1.5ML 3 AMP.") which now occupies row 10, pushing
# "NORHINOSE ..." down to row 11. Row 12 ("سرنجات 3 سم") is unchanged. The
# grand total in P13 and the generation timestamp in A14 are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: becomes the new "MOBITIL" item -----------------------------
$ws.Range("C10").Value = "MOBITIL 15MG/1.5ML 3 AMP."
$ws.Range("H10").Value = "2:2"

$fmt = $ws.Range("N10").NumberFormat
$ws.Range("N10").NumberFormat = "@"
$ws.Range("N10").Value = "39.00"
$ws.Range("N10").NumberFormat = $fmt

$fmt = $ws.Range("P10").NumberFormat
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "25.7400"
$ws.Range("P10").NumberFormat = $fmt

$ws.Range("Q10").Value = "0:2"

# --- Row 11: now holds the "NORHINOSE" item (shifted down from row 10) ---
$ws.Range("C11").Value = "NORHINOSE 50MCG/DOSE NASAL SPRAY 120 DOSES"
$ws.Range("H11").Value = "3:0"

$fmt = $ws.Range("N11").NumberFormat
$ws.Range("N11").NumberFormat = "@"
$ws.Range("N11").Value = "90.00"
$ws.Range("N11").NumberFormat = $fmt

$fmt = $ws.Range("P11").NumberFormat
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "90.0000"
$ws.Range("P11").NumberFormat = $fmt

$ws.Range("Q11").Value = "1:0"

# --- Row 12 ("سرنجات 3 سم") is unchanged, left as-is ----------------------

# --- Grand total -----------------------------------------------------------
$ws.Range("P13").Value = 287.74

# --- Footer: regenerate timestamp ------------------------------------------
$ws.Range("A14").Value = "Friday, 15 August, 2025 4:19 PM"
